# Finalized definition of virtual environment
# Updates the calibration results (param values + abs_epsi_autocorr arrays)
# on each of the risk-driver-type model sheets.

$wb = $excel.ActiveWorkbook

# --- AR sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.1776109831868101
$ws.Range("B3").Value = 0.7783110693441559
$ws.Range("B4").Value = 383.848296712782
$ws.Range("B5").Value = "[1.0, 0.1810680655162657, 0.1389408881107738, 0.13568495174960218, 0.17438166284207746, 0.3968835663359785, 0.14958086177894703, 0.1972277462259116, 0.17187336246720356, 0.1647020588388224, 0.12319251699959385, 0.14287735740216964, 0.13409210900116061, 0.1896551839009485, 0.2305464903146036, 0.13379265951828337, 0.1692363483153363, 0.13682927450674534, 0.12145623928074482, 0.13287253682205427]"

# --- SETAR sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value = 1.06969696969697
$ws.Range("B4").Value = -7.026667840185721
$ws.Range("B5").Value = 0.6879047765065119
$ws.Range("B6").Value = 249.4505626956585
$ws.Range("B7").Value = 6.864561894726315
$ws.Range("B8").Value = 0.6922773217549091
$ws.Range("B9").Value = 242.8517483897935
$ws.Range("B10").Value = "[0.9999999999999999, 0.24771717609180904, 0.1829897222994545, 0.16179815961692046, 0.20801240975489613, 0.2865171501865134, 0.12232681364711173, 0.15057458856620715, 0.13209976169096183, 0.16842724292567463, 0.17231530964410674, 0.15077633018114048, 0.11854940902901025, 0.15737442471335736, 0.1818699312722914, 0.12559955900207448, 0.1550310318897615, 0.16813727829129818, 0.14142357826173899, 0.11102566554428152]"

# --- GARCH sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.003620499177789568
$ws.Range("B3").Value = 3.556840396163731
$ws.Range("B4").Value = 0.03985447404472094
$ws.Range("B5").Value = 0.9488371413387855
$ws.Range("B6").Value = "[1.0, 0.013664902597682424, -0.04835544746087865, -0.07415039606176063, -0.06372547352477723, 0.30140981848791903, 0.02029249152402429, 0.009248824276025113, 0.015127575086847032, -0.014827632068282373, -0.08776497467455034, -0.059657727511022794, -0.028542689413831337, 0.027979236811910373, 0.04732719241902221, -0.01568598381192972, -0.010146843157736795, -0.03770774039442543, -0.04462489898131445, -0.038832468133576215]"

# --- TARCH sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.02644364686685877
$ws.Range("B3").Value = 3.562689117399171
$ws.Range("B4").Value = 0.0382910043017482
$ws.Range("B5").Value = 0.9487401072585366
$ws.Range("B6").Value = "[0.9999999999999998, 0.013431549894932135, -0.04839480551473856, -0.07429445780478955, -0.06373787905723027, 0.30154552974757, 0.02029286424551268, 0.009306580916913687, 0.014928720855878973, -0.014765373504382914, -0.08750042905560264, -0.05951887557354756, -0.028456387001959467, 0.027935617448488433, 0.047494530975396174, -0.0155871904385998, -0.009940750808888215, -0.03758606500052652, -0.04452323584282263, -0.038831449357241525]"
$ws.Range("B7").Value = 0.003313479378135294

# --- AR_TARCH sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.1376636876201886
$ws.Range("B3").Value = 3.279054323874111
$ws.Range("B4").Value = 0.03104880442516395
$ws.Range("B5").Value = 0.9526466683261864
$ws.Range("B6").Value = "[1.0, 0.024199944136474874, -0.04071684682780907, -0.06649319733123656, -0.021191906064247243, 0.23881938379890635, 0.013698181430104112, 0.01560507796895629, -0.011957705589616518, -0.018363886056220313, -0.0677658673827169, -0.038540301036076634, -0.022440500460462, 0.014577047628225313, 0.03268765136299908, -0.019231333787521843, -0.010021672011287956, -0.027497552047707374, -0.045878036567388966, -0.03512832346631533]"
$ws.Range("B7").Value = 0.0088581106396615
$ws.Range("B9").Value = 0.7810757852440677
